$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Register the small (size 8) font used for the phonetic-info settings,
# without leaving any visible residue on the grid (format a scratch cell,
# then clear it completely).
$ws.Range("Z100").Font.Size = 8
$ws.Range("Z100").Clear()

# --- Give column B the same look as column A (copy formats only, values are
# set further below) before typing the new data so the second column reuses
# the existing header / body cell styles instead of creating new ones.
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B4").PasteSpecial(-4122)

# --- Fetch value from 2 columns: fill in the new product/credential data,
# in row-major order so shared strings get (re)built in a predictable order.
$ws.Range("B1").Value = "product"
$ws.Range("A2").Value = "student"
$ws.Range("B2").Value = "Password123"
$ws.Range("A3").Value = "student"
$ws.Range("B3").Value = "Password123"
$ws.Range("A4").Value = "student1"
$ws.Range("B4").Value = "Password123"

# --- Extend report config: widen column B and move the active selection.
$ws.Columns("B").ColumnWidth = 13.15
[void]$ws.Range("A7").Select()
